$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2191780821917808
$ws.Range("C2").Value = 0.5178082191780822
$ws.Range("J2").Value = 0.0136986301369863
$ws.Range("P2").Value = 0.1506849315068493
$ws.Range("S2").Value = 0.09863013698630137
$ws.Range("B3").Value = 0.02010050251256281
$ws.Range("C3").Value = 0.03015075376884422
$ws.Range("J3").Value = 0.01507537688442211
$ws.Range("P3").Value = 0.7185929648241206
$ws.Range("S3").Value = 0.2160804020100502
$ws.Range("J4").Value = 0.09803921568627451
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.2352941176470588
$ws.Range("P5").Value = 1
$ws.Range("B6").Value = 0.0797872340425532
$ws.Range("D6").Value = 0.01595744680851064
$ws.Range("F6").Value = 0.06914893617021277
$ws.Range("J6").Value = 0.2340425531914894
$ws.Range("O6").Value = 0.02127659574468085
$ws.Range("Q6").Value = 0.1595744680851064
$ws.Range("R6").Value = 0.05319148936170213
$ws.Range("S6").Value = 0.3670212765957447
$ws.Range("B7").Value = 0.1421800947867299
$ws.Range("D7").Value = 0.02843601895734597
$ws.Range("F7").Value = 0.02843601895734597
$ws.Range("J7").Value = 0.1469194312796208
$ws.Range("O7").Value = 0.004739336492890996
$ws.Range("Q7").Value = 0.1895734597156398
$ws.Range("R7").Value = 0.05213270142180094
$ws.Range("S7").Value = 0.4075829383886256
$ws.Range("B8").Value = 0.09515260323159784
$ws.Range("D8").Value = 0.02154398563734291
$ws.Range("F8").Value = 0.04308797127468582
$ws.Range("J8").Value = 0.1310592459605027
$ws.Range("O8").Value = 0.0125673249551167
$ws.Range("Q8").Value = 0.1633752244165171
$ws.Range("R8").Value = 0.1149012567324955
$ws.Range("S8").Value = 0.4183123877917415
$ws.Range("B9").Value = 0.07954545454545454
$ws.Range("D9").Value = 0.0303030303030303
$ws.Range("F9").Value = 0.04924242424242424
$ws.Range("J9").Value = 0.1212121212121212
$ws.Range("O9").Value = 0.01136363636363636
$ws.Range("Q9").Value = 0.1856060606060606
$ws.Range("R9").Value = 0.1136363636363636
$ws.Range("S9").Value = 0.4090909090909091
$ws.Range("B10").Value = 0.1135714285714286
$ws.Range("D10").Value = 0.01714285714285714
$ws.Range("E10").Value = 0.0007142857142857143
$ws.Range("F10").Value = 0.05357142857142857
$ws.Range("J10").Value = 0.1314285714285714
$ws.Range("O10").Value = 0.01142857142857143
$ws.Range("Q10").Value = 0.2
$ws.Range("R10").Value = 0.08571428571428572
$ws.Range("S10").Value = 0.3864285714285715
$ws.Range("G11").Value = 0.1498257839721254
$ws.Range("J11").Value = 0.06968641114982578
$ws.Range("K11").Value = 0.1846689895470383
$ws.Range("L11").Value = 0.5853658536585366
$ws.Range("S11").Value = 0.01045296167247387
$ws.Range("G12").Value = 0.7430167597765364
$ws.Range("J12").Value = 0.1731843575418995
$ws.Range("L12").Value = 0.05027932960893855
$ws.Range("S12").Value = 0.0335195530726257
$ws.Range("G13").Value = 0.7454545454545455
$ws.Range("J13").Value = 0.2
$ws.Range("S13").Value = 0.05454545454545454
$ws.Range("H15").Value = 0.1674008810572687
$ws.Range("I15").Value = 0.1013215859030837
$ws.Range("J15").Value = 0.3612334801762114
$ws.Range("K15").Value = 0.08370044052863436
$ws.Range("M15").Value = 0.01762114537444934
$ws.Range("O15").Value = 0.07488986784140969
$ws.Range("S15").Value = 0.1938325991189427
$ws.Range("F16").Value = 0.008968609865470852
$ws.Range("H16").Value = 0.1569506726457399
$ws.Range("I16").Value = 0.1121076233183857
$ws.Range("J16").Value = 0.3991031390134529
$ws.Range("K16").Value = 0.1076233183856502
$ws.Range("M16").Value = 0.04035874439461883
$ws.Range("O16").Value = 0.04932735426008968
$ws.Range("S16").Value = 0.1255605381165919
$ws.Range("F17").Value = 0.02053388090349076
$ws.Range("H17").Value = 0.1930184804928131
$ws.Range("I17").Value = 0.1273100616016427
$ws.Range("J17").Value = 0.4291581108829569
$ws.Range("K17").Value = 0.06776180698151951
$ws.Range("M17").Value = 0.02258726899383984
$ws.Range("O17").Value = 0.04928131416837783
$ws.Range("S17").Value = 0.09034907597535935
$ws.Range("F18").Value = 0.0170940170940171
$ws.Range("H18").Value = 0.2136752136752137
$ws.Range("I18").Value = 0.1111111111111111
$ws.Range("J18").Value = 0.4401709401709402
$ws.Range("K18").Value = 0.03418803418803419
$ws.Range("M18").Value = 0.0170940170940171
$ws.Range("N18").Value = 0.004273504273504274
$ws.Range("O18").Value = 0.05128205128205128
$ws.Range("S18").Value = 0.1111111111111111
$ws.Range("F19").Value = 0.01403508771929825
$ws.Range("H19").Value = 0.2385964912280702
$ws.Range("I19").Value = 0.09333333333333334
$ws.Range("J19").Value = 0.3466666666666667
$ws.Range("K19").Value = 0.1017543859649123
$ws.Range("M19").Value = 0.01964912280701754
$ws.Range("N19").Value = 0.001403508771929824
$ws.Range("O19").Value = 0.07157894736842105
$ws.Range("S19").Value = 0.1129824561403509
